$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the label text in A2
$ws.Range("A2").Value = "1dq ED 01 a1"

# Move the active selection to A3 (matches the saved selection in the diff)
$ws.Range("A3").Select()
